$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -9
$ws.Range("F6").Value = -10
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -9
